$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-annotated DAMSLTag (col I) / DialogAct (col J) pairs per row, per SGNN re-run
$ws.Range("I38").Value = "ba"
$ws.Range("J38").Value = "Appreciation"
$ws.Range("I39").Value = "sd"
$ws.Range("J39").Value = "Statement-non-opinion"
$ws.Range("I54").Value = "sv"
$ws.Range("J54").Value = "Statement-opinion"
$ws.Range("I57").Value = "sd"
$ws.Range("J57").Value = "Statement-non-opinion"
$ws.Range("I58").Value = "sd"
$ws.Range("J58").Value = "Statement-non-opinion"
$ws.Range("I67").Value = "ba"
$ws.Range("J67").Value = "Appreciation"
$ws.Range("I80").Value = "sd"
$ws.Range("J80").Value = "Statement-non-opinion"
$ws.Range("I81").Value = "sd"
$ws.Range("J81").Value = "Statement-non-opinion"
$ws.Range("I82").Value = "sv"
$ws.Range("J82").Value = "Statement-opinion"
$ws.Range("I91").Value = "sv"
$ws.Range("J91").Value = "Statement-opinion"
$ws.Range("I92").Value = "sd"
$ws.Range("J92").Value = "Statement-non-opinion"
$ws.Range("I102").Value = "aa"
$ws.Range("J102").Value = "Agree/Accept"
$ws.Range("I118").Value = "sd"
$ws.Range("J118").Value = "Statement-non-opinion"
$ws.Range("I119").Value = "sd"
$ws.Range("J119").Value = "Statement-non-opinion"
$ws.Range("I126").Value = "aa"
$ws.Range("J126").Value = "Agree/Accept"
$ws.Range("I127").Value = "aa"
$ws.Range("J127").Value = "Agree/Accept"
$ws.Range("I136").Value = "aa"
$ws.Range("J136").Value = "Agree/Accept"
$ws.Range("I178").Value = "ba"
$ws.Range("J178").Value = "Appreciation"
$ws.Range("I181").Value = "ba"
$ws.Range("J181").Value = "Appreciation"
$ws.Range("I182").Value = "b"
$ws.Range("J182").Value = "Acknowledge (Backchannel)"
$ws.Range("I185").Value = "sv"
$ws.Range("J185").Value = "Statement-opinion"
$ws.Range("I192").Value = "aa"
$ws.Range("J192").Value = "Agree/Accept"
$ws.Range("I194").Value = "aa"
$ws.Range("J194").Value = "Agree/Accept"
$ws.Range("I196").Value = "aa"
$ws.Range("J196").Value = "Agree/Accept"
$ws.Range("I202").Value = "aa"
$ws.Range("J202").Value = "Agree/Accept"
$ws.Range("I203").Value = "sd"
$ws.Range("J203").Value = "Statement-non-opinion"
$ws.Range("I206").Value = "aa"
$ws.Range("J206").Value = "Agree/Accept"
$ws.Range("I240").Value = "aa"
$ws.Range("J240").Value = "Agree/Accept"
$ws.Range("I242").Value = "%"
$ws.Range("J242").Value = "Uninterpretable"
$ws.Range("I244").Value = "aa"
$ws.Range("J244").Value = "Agree/Accept"
$ws.Range("I262").Value = "aa"
$ws.Range("J262").Value = "Agree/Accept"
$ws.Range("I264").Value = "aa"
$ws.Range("J264").Value = "Agree/Accept"
$ws.Range("I280").Value = "aa"
$ws.Range("J280").Value = "Agree/Accept"
$ws.Range("I312").Value = "sd"
$ws.Range("J312").Value = "Statement-non-opinion"
$ws.Range("I314").Value = "aa"
$ws.Range("J314").Value = "Agree/Accept"
$ws.Range("I316").Value = "ba"
$ws.Range("J316").Value = "Appreciation"
$ws.Range("I320").Value = "sd"
$ws.Range("J320").Value = "Statement-non-opinion"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I350").Value = "aa"
$ws.Range("J350").Value = "Agree/Accept"
